$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80, shifting existing rows 80..169 down to 81..170
$ws.Rows.Item(80).Insert()

# Populate the new row 80 with the new weekly price record
$ws.Range("A80").Value() = 10
$ws.Range("B80").Value() = "Vega Modelo de Temuco"
$ws.Range("C80").Value() = "La Araucanía"
$ws.Range("D80").Value() = 44579
$ws.Range("E80").Value() = 9
$ws.Range("F80").Value() = 100112043
$ws.Range("G80").Value() = "Pepino dulce"
$ws.Range("H80").Value() = "Cultivar IV Región"
$ws.Range("I80").Value() = "Primera"
$ws.Range("J80").Value() = 40
$ws.Range("K80").Value() = 24000
$ws.Range("L80").Value() = 24000
$ws.Range("M80").Value() = 24000
$ws.Range("N80").Value() = "`$/bandeja 18 kilos"
$ws.Range("O80").Value() = "Provincia de Limarí"
$ws.Range("P80").Value() = 1333
$ws.Range("Q80").Value() = 18
$ws.Range("R80").Value() = "Hortaliza"

Write-Host ("Done. UsedRange rows: " + $ws.UsedRange.Rows.Count)
